$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new column before column A; every existing column (and its data,
# formatting, etc.) shifts one position to the right (B->C, C->D, ...).
$ws.Columns.Item(1).Insert()

# Set the header text for the brand-new column A.
$ws.Range("A1").Value = "客户主体"

# Give the new header cell A1 the same look (bold, bordered, centered) as
# the rest of row 1 by copying the formatting from the neighboring header
# cell (now in column B).
$ws.Range("B1").Copy() | Out-Null
$ws.Range("A1").PasteSpecial(-4122) | Out-Null

# Fill in the auto-detected client type ("enterprise" for companies,
# "personal" for individuals) for each data row in the new column A.
$ws.Range("A2").Value = "personal"
$ws.Range("A3").Value = "personal"
$ws.Range("A4").Value = "enterprise"
$ws.Range("A5").Value = "personal"
$ws.Range("A6").Value = "enterprise"
$ws.Range("A7").Value = "personal"
